$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-05-02 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-03 Friday", 1) | Out-Null

$tbl = $d.Tables.Item(1)

# Row 1
$cell = $tbl.Cell(1, 1)
$cell.Range.Find.Execute("92÷4=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=6, 0", 1) | Out-Null
$cell = $tbl.Cell(1, 2)
$cell.Range.Find.Execute("77÷5=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=6, 0", 1) | Out-Null
$cell = $tbl.Cell(1, 3)
$cell.Range.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷8=2, 3", 1) | Out-Null
$cell = $tbl.Cell(1, 4)
$cell.Range.Find.Execute("53÷7=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "33÷7=4, 5", 1) | Out-Null
$cell = $tbl.Cell(1, 5)
$cell.Range.Find.Execute("61÷2=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "25÷5=5, 0", 1) | Out-Null

# Row 5
$cell = $tbl.Cell(5, 1)
$cell.Range.Find.Execute("82÷4=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "89÷5=17, 4", 1) | Out-Null
$cell = $tbl.Cell(5, 2)
$cell.Range.Find.Execute("78÷7=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "23÷6=3, 5", 1) | Out-Null
$cell = $tbl.Cell(5, 3)
$cell.Range.Find.Execute("18÷5=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "52÷4=13, 0", 1) | Out-Null
$cell = $tbl.Cell(5, 4)
$cell.Range.Find.Execute("85÷7=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "55÷8=6, 7", 1) | Out-Null
$cell = $tbl.Cell(5, 5)
$cell.Range.Find.Execute("15÷9=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=42, 0", 1) | Out-Null

# Row 9
$cell = $tbl.Cell(9, 1)
$cell.Range.Find.Execute("77÷9=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=1, 7", 1) | Out-Null
$cell = $tbl.Cell(9, 2)
$cell.Range.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "61÷6=10, 1", 1) | Out-Null
$cell = $tbl.Cell(9, 3)
$cell.Range.Find.Execute("60÷4=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "26÷2=13, 0", 1) | Out-Null
$cell = $tbl.Cell(9, 4)
$cell.Range.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "72÷4=18, 0", 1) | Out-Null
$cell = $tbl.Cell(9, 5)
$cell.Range.Find.Execute("78÷2=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "16÷5=3, 1", 1) | Out-Null

# Row 13
$cell = $tbl.Cell(13, 1)
$cell.Range.Find.Execute("54÷2=27, 0", $true, $false, $false, $false, $false, $true, 1, $false, "91÷6=15, 1", 1) | Out-Null
$cell = $tbl.Cell(13, 2)
$cell.Range.Find.Execute("53÷2=26, 1", $true, $false, $false, $false, $false, $true, 1, $false, "79÷7=11, 2", 1) | Out-Null
$cell = $tbl.Cell(13, 3)
$cell.Range.Find.Execute("67÷5=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "93÷4=23, 1", 1) | Out-Null
$cell = $tbl.Cell(13, 4)
$cell.Range.Find.Execute("22÷6=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "33÷5=6, 3", 1) | Out-Null
$cell = $tbl.Cell(13, 5)
$cell.Range.Find.Execute("74÷8=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "90÷5=18, 0", 1) | Out-Null

# Row 17
$cell = $tbl.Cell(17, 1)
$cell.Range.Find.Execute("67÷3=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=8, 4", 1) | Out-Null
$cell = $tbl.Cell(17, 2)
$cell.Range.Find.Execute("57÷7=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷4=10, 0", 1) | Out-Null
$cell = $tbl.Cell(17, 3)
$cell.Range.Find.Execute("16÷4=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "81÷5=16, 1", 1) | Out-Null
$cell = $tbl.Cell(17, 4)
$cell.Range.Find.Execute("97÷5=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "22÷3=7, 1", 1) | Out-Null
$cell = $tbl.Cell(17, 5)
$cell.Range.Find.Execute("20÷4=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "27÷3=9, 0", 1) | Out-Null
